# Append the 14/04/2019 .. 27/04/2019 logbook entries to the end of the
# document. The existing "_GoBack" bookmark (which originally sat right
# after the 12/04/2019 entry's final sentence) is removed and re-inserted
# in its new location -- mid-sentence inside the final (27/04/2019) entry,
# exactly where the author's cursor last was -- as part of the inserted
# XML fragment below.

$d = $word.ActiveDocument

# Drop the pre-existing "_GoBack" bookmark; it is recreated below inside
# the newly inserted content, in its new resting place.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete() | Out-Null

# Collapse a range onto the very end of the document body (immediately
# after "...golang." which closes out the 12/04/2019 entry) and splice in
# all of the new paragraphs as raw WordprocessingML.
$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos, $endPos)

$newEntriesXml = @'
<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>14/04/2019</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Go </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>lang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> http handler did not work, so I decided to copy with full credit to the original http handler, but add a line where I pass the variable within the context.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>15/04/2019</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">I started working on the frontend, using react router and react to create prototypes of the login screen and the home screen using material </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>16/04/2019</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Another set of refactoring because I realised that some of the functions of the server side, namely the resolver to get a login token </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>could be shrunk</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">I had a set of three functions to handle each user type, but it was three copies of the same function, and in an effort to make the application DRY, I created one function that handles the signing in of all users, instead to login in different user types the application now has an </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>enum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> where you can select user types with a select within the frontend html select statement.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> This </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>gets</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> passed as an argument to the resolver function which uses a case statement to login a user of all different types. The problem is when the form submits, it passes all the variables through the html submit functionality. So I started working on react state management as opposed to the traditional html submit to handle the data.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>17/04/2019</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Ditched the idea of using a switch to handle the user selection for logging in because I realised it </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>wasn’t</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the best way of doing things. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Instead</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> I use the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, and have an argument </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>usertype</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> which </w:t></w:r><w:r><w:t xml:space="preserve">passes what type of login the user wants. The argument </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>is saved</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> as a state.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>19</w:t></w:r><w:r><w:t>/04/2019</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Finished the form completely for login, with routing and handling </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>graphql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> errors. The only problem is the application is returning </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> overflow error and I </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>cant</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> find the source of the bug. I see to many state calls from the login form, but that could </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>either be</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the router, the login form or the Apollo library complaining. Will try to fix it but in the meantime improving the backend with more refactoring of functions. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>24/04/2019</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Finished all the authentication and authorization on the backend, and I have started work on fixing the bug. I could not find what was causing the bug, and I decided rather than trying to fix the bug, I would rewrite the entire frontend, using typescript for type checking and making it purely functional. That includes avoiding using loops and instead using iterator methods and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>callback</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> functions, using constants for everything. I use react hooks for state management so I can have functional forms instead of complicated classes with a single method, and have </w:t></w:r><w:r><w:t xml:space="preserve">removed saving the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>usertype</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as a state because the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>usertype</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> never changes on a particular </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">On the server side, I have done the complete opposite and have refactored the code to better suit an </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>object oriented</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> paradigm. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>25/04/2019</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">I refactored the frontend and improved my typescript, before I used to pass the props as a type of any, but this made the entire type checking redundant. So I </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>learn’t</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> how to create interfaces which are used to add type checking to certain functions, I have switched from Apollo to Apollo-hooks because it brings down the amount of code I need to write and simplifies a lot of things.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>26/04/2019</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Today was testing the frontend and backend together and completely finishing authentication and authorization. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>27/04/2019</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Today was fixing bugs in the backend and developing the actual core functionality. I completely removed the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>env.json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> idea and now I hard code the variables into a go file, I added a commented line in the .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gitignore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file that ignores a go file that has these hard coded variables. I included this in the build so that the project could compile online. All other variables </w:t></w:r><w:r><w:t xml:space="preserve">apart from the secret key will be stored in the database. I added more object oriented style programming to the backend </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>through the use of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> getter methods of the environment variables and also separated the database from the schema through a private pointer and a getter method to the database connection. A lot of the code for the secret key </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>was removed</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">, now it is auto generated </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>everytime</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the application restarts and is stored in a temporary variable. This means that every time the server restarts, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>everybodies</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tokens will be </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>invaled</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and they will need to sign in again. Through </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>this</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> there is less code and the application is more secure, if I have time I might even make the tokens secret key</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> regenerate each week increasing the security even more.</w:t></w:r></w:p>
'@

$insertionPoint.InsertXML($newEntriesXml) | Out-Null
